$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Brenselsstoffer"
$ws.Range("A9").Value = "Kjemiske produkter"

$ws.Range("F9").Select()
